$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the date-sorted Camote block (rows 329-330),
# pushing the existing rows 329:373 down to 331:375.
$ws.Rows("329:330").Insert()

# Row 329: new weekly entry - "1a nueva(o)" from Peru
$ws.Range("A329").Value = 8
$ws.Range("B329").Value = "Terminal La Palmera de La Serena"
$ws.Range("C329").Value = "Coquimbo"
$ws.Range("D329").Value = 44474
$ws.Range("E329").Value = 4
$ws.Range("F329").Value = 100112045
$ws.Range("G329").Value = "Zapallo"
$ws.Range("H329").Value = "Camote"
$ws.Range("I329").Value = "1a nueva(o)"
$ws.Range("J329").Value = 800
$ws.Range("K329").Value = 700
$ws.Range("L329").Value = 750
$ws.Range("M329").Value = 725
$ws.Range("N329").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O329").Value = "Perú"
$ws.Range("P329").Value = 725
$ws.Range("Q329").Value = 1
$ws.Range("R329").Value = "Hortaliza"

# Row 330: new weekly entry - "2a nueva(o)" from Peru
$ws.Range("A330").Value = 8
$ws.Range("B330").Value = "Terminal La Palmera de La Serena"
$ws.Range("C330").Value = "Coquimbo"
$ws.Range("D330").Value = 44474
$ws.Range("E330").Value = 4
$ws.Range("F330").Value = 100112045
$ws.Range("G330").Value = "Zapallo"
$ws.Range("H330").Value = "Camote"
$ws.Range("I330").Value = "2a nueva(o)"
$ws.Range("J330").Value = 540
$ws.Range("K330").Value = 600
$ws.Range("L330").Value = 650
$ws.Range("M330").Value = 625
$ws.Range("N330").Value = "`$/kilo (volumen en unidades)"
$ws.Range("O330").Value = "Perú"
$ws.Range("P330").Value = 625
$ws.Range("Q330").Value = 1
$ws.Range("R330").Value = "Hortaliza"
